$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Column G width change: 17 -> 18 characters
$ws1.Columns.Item(7).ColumnWidth = 17.14

# Row 3
$ws1.Range("G3").Value = "Halftime"

# Row 4
$ws1.Range("G4").Value = "18:10 - 2nd Half"
$ws1.Range("H4").Value = 4
$ws1.Range("J4").Value = 2
$ws1.Range("O4").Value = 14

# Row 7
$ws1.Range("G7").Value = "18:10 - 2nd Half"
$ws1.Range("H7").Value = 14
$ws1.Range("J7").Value = 4
$ws1.Range("O7").Value = 18

# Row 9
$ws1.Range("G9").Value = "18:10 - 2nd Half"

# Row 10
$ws1.Range("G10").Value = "Halftime"

# Row 11
$ws1.Range("G11").Value = "18:10 - 2nd Half"
$ws1.Range("H11").Value = -4
$ws1.Range("O11").Value = 14

# Row 12
$ws1.Range("G12").Value = "Halftime"

# Row 13
$ws1.Range("G13").Value = "Halftime"

# Row 14
$ws1.Range("G14").Value = "Halftime"

# Row 15
$ws1.Range("G15").Value = "18:10 - 2nd Half"
$ws1.Range("H15").Value = 9
$ws1.Range("J15").Value = 4
$ws1.Range("N15").Value = 1
$ws1.Range("O15").Value = 19

# Row 17
$ws1.Range("G17").Value = "18:10 - 2nd Half"
$ws1.Range("H17").Value = 9
$ws1.Range("I17").Value = 9
$ws1.Range("O17").Value = 18

# Row 18
$ws1.Range("G18").Value = "Halftime"

# Row 19
$ws1.Range("G19").Value = "18:10 - 2nd Half"

# Row 20
$ws1.Range("G20").Value = "Halftime"

# Row 21
$ws1.Range("D21").Value = "Xzayvier Brown"
$ws1.Range("E21").Value = "OU"
$ws1.Range("G21").Value = "18:10 - 2nd Half"
$ws1.Range("H21").Value = 8
$ws1.Range("I21").Value = 11
$ws1.Range("J21").Value = 1
$ws1.Range("K21").Value = 0
$ws1.Range("L21").Value = 0
$ws1.Range("O21").Value = 21

# Row 22
$ws1.Range("D22").Value = "Xaivian Lee"
$ws1.Range("E22").Value = "FLA"
$ws1.Range("G22").Value = "18:10 - 2nd Half"
$ws1.Range("H22").Value = 6
$ws1.Range("I22").Value = 4
$ws1.Range("J22").Value = 2
$ws1.Range("K22").Value = 2
$ws1.Range("L22").Value = 1
$ws1.Range("O22").Value = 12

# Row 23
$ws1.Range("G23").Value = "18:10 - 2nd Half"

# Row 25
$ws1.Range("G25").Value = "18:10 - 2nd Half"
$ws1.Range("H25").Value = 13
$ws1.Range("I25").Value = 7
$ws1.Range("O25").Value = 15

# Row 26
$ws1.Range("G26").Value = "18:10 - 2nd Half"
$ws1.Range("H26").Value = 6
$ws1.Range("K26").Value = 1
$ws1.Range("O26").Value = 16

# Row 27
$ws1.Range("G27").Value = "Halftime"

# Row 29
$ws1.Range("G29").Value = "18:10 - 2nd Half"

# Row 31
$ws1.Range("G31").Value = "Halftime"

# Row 33
$ws1.Range("G33").Value = "Halftime"

# Row 34
$ws1.Range("G34").Value = "18:10 - 2nd Half"
$ws1.Range("H34").Value = 5
$ws1.Range("L34").Value = 1
$ws1.Range("O34").Value = 19

# Row 38
$ws1.Range("G38").Value = "Halftime"

# Row 40
$ws1.Range("G40").Value = "Halftime"

# Row 41
$ws1.Range("G41").Value = "18:10 - 2nd Half"

# Row 43
$ws1.Range("G43").Value = "Halftime"

# Row 44
$ws1.Range("G44").Value = "Halftime"

# Row 48
$ws1.Range("G48").Value = "18:10 - 2nd Half"

# Row 49
$ws1.Range("G49").Value = "Halftime"

# Row 54
$ws1.Range("G54").Value = "18:10 - 2nd Half"

# Row 55
$ws1.Range("G55").Value = "Halftime"

# Row 56
$ws1.Range("G56").Value = "18:10 - 2nd Half"

# Row 57
$ws1.Range("G57").Value = "Halftime"

# OwnerTotals sheet updates
$ws2.Range("B2").Value = 36
$ws2.Range("B3").Value = 28
$ws2.Range("B5").Value = 20
$ws2.Range("B6").Value = 9
